$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 111
$ws.Range("H111").Value = 2836.6667
$ws.Range("I111").Value = 1397.6
$ws.Range("J111").Value = 10032
$ws.Range("K111").Value = 4192.799999999999
$ws.Range("L111").Value = 30096
$ws.Range("M111").Value = -1125.799999999999
$ws.Range("N111").Value = -36230

# Row 131
$ws.Range("H131").Value = 3698.75
$ws.Range("I131").Value = 1795
$ws.Range("J131").Value = 4333.3335
$ws.Range("K131").Value = 5385
$ws.Range("L131").Value = 13000.0005
$ws.Range("M131").Value = -345
$ws.Range("N131").Value = -23080.0005

# Row 137
$ws.Range("H137").Value = 951.5454999999999
$ws.Range("I137").Value = 958.8
$ws.Range("K137").Value = 2876.4
$ws.Range("M137").Value = -326.3999999999996

# Row 138
$ws.Range("H138").Value = 2300.652
$ws.Range("J138").Value = 2348.2273
$ws.Range("L138").Value = 7044.6819
$ws.Range("N138").Value = -17324.6819

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 15000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 15000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 15000
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -15346

# Row 12
$ws.Range("H12").Value = 25500
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 50000
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 50000
$ws.Range("M12").Value = -827
$ws.Range("N12").Value = -50346

# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Row 95
$ws.Range("H95").Value = 24402.666
$ws.Range("J95").Value = 24402.666
$ws.Range("L95").Value = 24402.666
$ws.Range("N95").Value = -29894.666

# Row 104
$ws.Range("H104").Value = 45000
$ws.Range("J104").Value = 45000
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -51988

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 18498.5
$ws.Range("J7").Value = 18498.5
$ws.Range("L7").Value = 18498.5
$ws.Range("N7").Value = -18724.5

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 134
$ws.Range("H134").Value = 2399.889
$ws.Range("I134").Value = 2349.875
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 7049.625
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -4514.625
$ws.Range("N134").Value = -13470

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1055.5
$ws.Range("I16").Value = 1499.5
$ws.Range("J16").Value = 611.5
$ws.Range("K16").Value = 1499.5
$ws.Range("L16").Value = 611.5
$ws.Range("M16").Value = -1212.5
$ws.Range("N16").Value = -1185.5

# Row 32
$ws.Range("H32").Value = 809
$ws.Range("I32").Value = 809
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 809
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -493
$ws.Range("N32").ClearContents()

# Row 86
$ws.Range("H86").Value = 5981.2856
$ws.Range("I86").Value = 8247.75
$ws.Range("J86").Value = 2959.3333
$ws.Range("K86").Value = 8247.75
$ws.Range("L86").Value = 2959.3333
$ws.Range("M86").Value = -7124.75
$ws.Range("N86").Value = -5205.3333

# Row 89
$ws.Range("H89").Value = 5981.2856
$ws.Range("I89").Value = 8247.75
$ws.Range("J89").Value = 2959.3333
$ws.Range("K89").Value = 41238.75
$ws.Range("L89").Value = 14796.6665
$ws.Range("M89").Value = -35622.75
$ws.Range("N89").Value = -26028.6665

# Row 94
$ws.Range("H94").Value = 1301.2222
$ws.Range("J94").Value = 1442.8
$ws.Range("L94").Value = 1442.8
$ws.Range("N94").Value = -2344.8

# Row 113
$ws.Range("H113").Value = 1055.5
$ws.Range("I113").Value = 1499.5
$ws.Range("J113").Value = 611.5
$ws.Range("K113").Value = 1499.5
$ws.Range("L113").Value = 611.5
$ws.Range("M113").Value = 670.5
$ws.Range("N113").Value = -4951.5

# Row 141
$ws.Range("H141").Value = 149980
$ws.Range("J141").Value = 149980
$ws.Range("L141").Value = 149980
$ws.Range("N141").Value = -160340

$ws = $wb.Worksheets.Item("CUL")
# Row 108
$ws.Range("H108").Value = 395
$ws.Range("I108").Value = 395
$ws.Range("K108").Value = 1185
$ws.Range("M108").Value = 1695

# Row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").ClearContents()

# Row 131
$ws.Range("H131").Value = 2926.0667
$ws.Range("I131").Value = 1243.3334
$ws.Range("J131").Value = 3113.037
$ws.Range("K131").Value = 3730.0002
$ws.Range("L131").Value = 9339.110999999999
$ws.Range("M131").Value = 1309.9998
$ws.Range("N131").Value = -19419.111

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 102
$ws.Range("H102").Value = 3912.2856
$ws.Range("I102").Value = 3083
$ws.Range("K102").Value = 3083
$ws.Range("M102").Value = -1461

# Row 132
$ws.Range("H132").Value = 4199.4443
$ws.Range("I132").Value = 3999.5
$ws.Range("K132").Value = 11998.5
$ws.Range("M132").Value = -9468.5

$ws = $wb.Worksheets.Item("LTW")
# Row 19
$ws.Range("H19").Value = 5999
$ws.Range("I19").Value = 5999
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 5999
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -5829
$ws.Range("N19").ClearContents()

# Row 40
$ws.Range("H40").Value = 3471.5715
$ws.Range("I40").Value = 3599.3333
$ws.Range("J40").Value = 2705
$ws.Range("K40").Value = 3599.3333
$ws.Range("L40").Value = 2705
$ws.Range("M40").Value = -3463.3333
$ws.Range("N40").Value = -2977

# Row 46
$ws.Range("H46").Value = 8500
$ws.Range("J46").Value = 26500
$ws.Range("L46").Value = 26500
$ws.Range("N46").Value = -26876

# Row 104
$ws.Range("H104").Value = 42249.75
$ws.Range("J104").Value = 42249.75
$ws.Range("L104").Value = 42249.75
$ws.Range("N104").Value = -49237.75

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 31548.8
$ws.Range("I54").Value = 30000
$ws.Range("J54").Value = 32581.334
$ws.Range("K54").Value = 30000
$ws.Range("L54").Value = 32581.334
$ws.Range("M54").Value = -29480
$ws.Range("N54").Value = -33621.334

# Row 136
$ws.Range("H136").Value = 5714
$ws.Range("I136").Value = 5714
$ws.Range("K136").Value = 17142
$ws.Range("M136").Value = -14592
